# Fix erroneous calculation values in column A (distance/result column) that
# had been computed incorrectly depending on direction.
# Commit message: "Hata giderme Yöne bağlı olarak hatalı hesaplama yapma sorunu giderildi."
# (Bugfix: incorrect calculation depending on direction has been fixed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    21  = 53.12
    107 = 90.93000000000001
    108 = 86.93000000000001
    109 = 82.94
    110 = 78.94
    111 = 74.95
    112 = 70.95
    113 = 66.95
    121 = 68.79000000000001
    122 = 68.56999999999999
    123 = 68.34
    124 = 60
    125 = 83.27
    126 = 106.54
    127 = 132.85
    128 = 122.7
    129 = 67.28
    132 = 68.84999999999999
    133 = 70.63
    134 = 68.5
    135 = 68.31
    136 = 54.72
    137 = 67.27
    138 = 67.26000000000001
    139 = 67.26000000000001
    140 = 67.25
    141 = 78.3
    142 = 43.07
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
